$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 801
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 67
$ws.Range("D2").Value = 65
$ws.Range("E2").Value = 52
$ws.Range("F2").Value = 45

# Row 3
$ws.Range("A3").Value = 1203
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 15
$ws.Range("D3").Value = 15
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 15

# Row 4
$ws.Range("A4").Value = 1001
$ws.Range("B4").Value = 18
$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 75
$ws.Range("E4").Value = 60
$ws.Range("F4").Value = 72

# Row 6
$ws.Range("A6").Value = 901
$ws.Range("B6").Value = 16
$ws.Range("C6").Value = 15
$ws.Range("D6").Value = 45
$ws.Range("E6").Value = 60
$ws.Range("F6").Value = 60

# Row 7
$ws.Range("A7").Value = 601
$ws.Range("B7").Value = 9
$ws.Range("C7").Value = 60
$ws.Range("D7").Value = 67
$ws.Range("E7").Value = 60
$ws.Range("F7").Value = 42

# Row 9
$ws.Range("A9").Value = 1201
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = 10

# Row 10
$ws.Range("A10").Value = 1202
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 10
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 10

# Row 11
$ws.Range("A11").Value = 101
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 15
$ws.Range("E11").Value = 60
$ws.Range("F11").Value = 15

# Row 12
$ws.Range("A12").Value = 902
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0

# Row 13
$ws.Range("A13").Value = 701
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 90
$ws.Range("D13").Value = 45
$ws.Range("E13").Value = 97
$ws.Range("F13").Value = 15

# Row 14
$ws.Range("A14").Value = 301
$ws.Range("B14").Value = 6
$ws.Range("C14").Value = 45
$ws.Range("D14").Value = 30
$ws.Range("E14").Value = 60
$ws.Range("F14").Value = 45

# Row 15
$ws.Range("A15").Value = 401
$ws.Range("B15").Value = 9
$ws.Range("C15").Value = 48
$ws.Range("D15").Value = 67
$ws.Range("E15").Value = 75
$ws.Range("F15").Value = 45

# Row 22
$ws.Range("A22").Value = 602
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 9

# Row 23
$ws.Range("A23").Value = 402
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0

